$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A34").Value = 111635445
$ws.Range("B34").Value = 89686
$ws.Range("E34").Value = 658
$ws.Range("F34").Value = "Rosenticka"
$ws.Range("G34").Value = "Rhodofomes roseus"
$ws.Range("H34").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q34").Value = 539972.1173992374
$ws.Range("R34").Value = 7198351.138093079
$ws.Range("Z34").Value = "09:30"
$ws.Range("AB34").Value = "09:30"
$ws.Range("AW34").Value = "Yasmine Kindlund"
$ws.Range("AX34").Value = "Yasmine Kindlund, Isak Vahlström"
$ws.Range("A35").Value = 111635444
$ws.Range("B35").Value = 89686
$ws.Range("D35").Value = "NT"
$ws.Range("E35").Value = 658
$ws.Range("F35").Value = "Rosenticka"
$ws.Range("G35").Value = "Rhodofomes roseus"
$ws.Range("H35").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q35").Value = 539849.8711390977
$ws.Range("R35").Value = 7198368.616090678
$ws.Range("Z35").Value = "09:56"
$ws.Range("AB35").Value = "09:56"
$ws.Range("AW35").Value = "Yasmine Kindlund"
$ws.Range("AX35").Value = "Yasmine Kindlund, Isak Vahlström"
$ws.Range("A36").Value = 111634866
$ws.Range("B36").Value = 77515
$ws.Range("E36").Value = 6425
$ws.Range("F36").Value = "Garnlav"
$ws.Range("G36").Value = "Alectoria sarmentosa"
$ws.Range("H36").Value = "(Ach.) Ach."
$ws.Range("Q36").Value = 539873.9909718054
$ws.Range("R36").Value = 7198345.158619706
$ws.Range("Z36").Value = "09:46"
$ws.Range("AB36").Value = "09:46"
$ws.Range("A37").Value = 111634867
$ws.Range("B37").Value = 73696
$ws.Range("E37").Value = 6440
$ws.Range("F37").Value = "Vitgrynig nållav"
$ws.Range("G37").Value = "Chaenotheca subroscida"
$ws.Range("H37").Value = "(Eitner) Zahlbr."
$ws.Range("Q37").Value = 539871.8034722162
$ws.Range("R37").Value = 7198349.800304586
$ws.Range("Z37").Value = "09:46"
$ws.Range("AB37").Value = "09:46"
$ws.Range("AW37").Value = "Isak Vahlström"
$ws.Range("AX37").Value = "Isak Vahlström, Yasmine Kindlund"
$ws.Range("A38").Value = 111635452
$ws.Range("B38").Value = 78578
$ws.Range("E38").Value = 6458
$ws.Range("F38").Value = "Lunglav"
$ws.Range("G38").Value = "Lobaria pulmonaria"
$ws.Range("H38").Value = "(L.) Hoffm."
$ws.Range("Q38").Value = 539861.2921981018
$ws.Range("R38").Value = 7198404.860384831
$ws.Range("Z38").Value = "09:50"
$ws.Range("AB38").Value = "09:50"
$ws.Range("A39").Value = 111635489
$ws.Range("B39").Value = 77515
$ws.Range("D39").Value = "NT"
$ws.Range("E39").Value = 6425
$ws.Range("F39").Value = "Garnlav"
$ws.Range("G39").Value = "Alectoria sarmentosa"
$ws.Range("H39").Value = "(Ach.) Ach."
$ws.Range("Q39").Value = 539945.9506927577
$ws.Range("R39").Value = 7198336.776317291
$ws.Range("Z39").Value = "09:39"
$ws.Range("AB39").Value = "09:39"
$ws.Range("A40").Value = 111635461
$ws.Range("B40").Value = 89590
$ws.Range("D40").Value = "VU"
$ws.Range("E40").Value = 48
$ws.Range("F40").Value = "Lappticka"
$ws.Range("G40").Value = "Amylocystis lapponica"
$ws.Range("H40").Value = "(Romell) Singer"
$ws.Range("Q40").Value = 539846.9353019162
$ws.Range("R40").Value = 7198365.604689348
$ws.Range("Z40").Value = "09:56"
$ws.Range("AB40").Value = "09:56"
$ws.Range("AW40").Value = "Yasmine Kindlund"
$ws.Range("AX40").Value = "Yasmine Kindlund, Isak Vahlström"
$ws.Range("A41").Value = 111634859
$ws.Range("B41").Value = 77515
$ws.Range("E41").Value = 6425
$ws.Range("F41").Value = "Garnlav"
$ws.Range("G41").Value = "Alectoria sarmentosa"
$ws.Range("H41").Value = "(Ach.) Ach."
$ws.Range("Q41").Value = 539847.161346367
$ws.Range("R41").Value = 7198348.622951495
$ws.Range("Z41").Value = "09:58"
$ws.Range("AB41").Value = "09:58"
$ws.Range("AW41").Value = "Isak Vahlström"
$ws.Range("AX41").Value = "Isak Vahlström, Yasmine Kindlund"
$ws.Range("A42").Value = 111634868
$ws.Range("B42").Value = 78612
$ws.Range("E42").Value = 6464
$ws.Range("F42").Value = "Luddlav"
$ws.Range("G42").Value = "Nephroma resupinatum"
$ws.Range("H42").Value = "(L.) Ach."
$ws.Range("Q42").Value = 539976.4302002029
$ws.Range("R42").Value = 7198378.371244119
$ws.Range("Z42").Value = "09:28"
$ws.Range("AB42").Value = "09:28"
$ws.Range("A43").Value = 111635413
$ws.Range("B43").Value = 89369
$ws.Range("D43").Value = "LC"
$ws.Range("E43").Value = 5447
$ws.Range("F43").Value = "Vedticka"
$ws.Range("G43").Value = "Fuscoporia viticola"
$ws.Range("H43").Value = "(Schwein.) Murrill"
$ws.Range("Q43").Value = 539850.8116781802
$ws.Range("R43").Value = 7198361.834730743
$ws.Range("Z43").Value = "09:57"
$ws.Range("AB43").Value = "09:57"
$ws.Range("AW43").Value = "Yasmine Kindlund"
$ws.Range("AX43").Value = "Yasmine Kindlund, Isak Vahlström"
$ws.Range("A44").Value = 111635499
$ws.Range("B44").Value = 85715
$ws.Range("D44").Value = "NT"
$ws.Range("E44").Value = 510
$ws.Range("F44").Value = "Doftskinn"
$ws.Range("G44").Value = "Cystostereum murrayi"
$ws.Range("H44").Value = "(Berk. & M.A. Curtis.) Pouzar"
$ws.Range("Q44").Value = 540009.9192712342
$ws.Range("R44").Value = 7198353.766191677
$ws.Range("Z44").Value = "09:32"
$ws.Range("AB44").Value = "09:32"
$ws.Range("A45").Value = 111635422
$ws.Range("B45").Value = 56398
$ws.Range("D45").Value = "NT"
$ws.Range("E45").Value = 100109
$ws.Range("F45").Value = "Tretåig hackspett"
$ws.Range("G45").Value = "Picoides tridactylus"
$ws.Range("H45").Value = "(Linnaeus, 1758)"
$ws.Range("M45").Value = "äldre spår"
$ws.Range("Q45").Value = 539953.4033757704
$ws.Range("R45").Value = 7198319.890847754
$ws.Range("Z45").Value = "09:36"
$ws.Range("AB45").Value = "09:36"
$ws.Range("A46").Value = 111635462
$ws.Range("B46").Value = 89590
$ws.Range("D46").Value = "VU"
$ws.Range("E46").Value = 48
$ws.Range("F46").Value = "Lappticka"
$ws.Range("G46").Value = "Amylocystis lapponica"
$ws.Range("H46").Value = "(Romell) Singer"
$ws.Range("Q46").Value = 539961.7289606878
$ws.Range("R46").Value = 7198365.011824355
$ws.Range("Z46").Value = "09:30"
$ws.Range("AB46").Value = "09:30"
$ws.Range("A47").Value = 111634869
$ws.Range("B47").Value = 78578
$ws.Range("E47").Value = 6458
$ws.Range("F47").Value = "Lunglav"
$ws.Range("G47").Value = "Lobaria pulmonaria"
$ws.Range("H47").Value = "(L.) Hoffm."
$ws.Range("Q47").Value = 539972.5933666634
$ws.Range("R47").Value = 7198379.169240371
$ws.Range("Z47").Value = "09:27"
$ws.Range("AB47").Value = "09:27"
$ws.Range("AW47").Value = "Isak Vahlström"
$ws.Range("AX47").Value = "Isak Vahlström, Yasmine Kindlund"
$ws.Range("A48").Value = 111634865
$ws.Range("B48").Value = 90087
$ws.Range("D48").Value = "LC"
$ws.Range("E48").Value = 3298
$ws.Range("F48").Value = "Trådticka"
$ws.Range("G48").Value = "Climacocystis borealis"
$ws.Range("H48").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q48").Value = 539879.8909062841
$ws.Range("R48").Value = 7198349.058794393
$ws.Range("Z48").Value = "09:47"
$ws.Range("AB48").Value = "09:47"
$ws.Range("AW48").Value = "Isak Vahlström"
$ws.Range("AX48").Value = "Isak Vahlström, Yasmine Kindlund"
$ws.Range("A49").Value = 111635437
$ws.Range("B49").Value = 89845
$ws.Range("D49").Value = "VU"
$ws.Range("E49").Value = 1209
$ws.Range("F49").Value = "Rynkskinn"
$ws.Range("G49").Value = "Phlebia centrifuga"
$ws.Range("H49").Value = "P.Karst."
$ws.Range("Q49").Value = 539973.573864806
$ws.Range("R49").Value = 7198369.416147546
$ws.Range("Z49").Value = "09:32"
$ws.Range("AB49").Value = "09:32"
$ws.Range("A50").Value = 111635419
$ws.Range("B50").Value = 89405
$ws.Range("E50").Value = 1202
$ws.Range("F50").Value = "Ullticka"
$ws.Range("G50").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H50").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q50").Value = 539844.8100177459
$ws.Range("R50").Value = 7198365.57640036
$ws.Range("Z50").Value = "09:58"
$ws.Range("AB50").Value = "09:58"
$ws.Range("M38").ClearContents()
